$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 1187
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 516.237
$ws.Range("E2").Value = 735.242
$ws.Range("H2").Value = 5795
$ws.Range("I2").Value = 19705
$ws.Range("J2").Value = 30
$ws.Range("K2").Value = 139
$ws.Range("L2").Value = 4617

$ws.Range("B3").Value = 263

$ws.Range("B4").Value = 179
$ws.Range("D4").Value = 58
$ws.Range("E4").Value = 75
$ws.Range("H4").Value = 2267
$ws.Range("I4").Value = 3103

$ws.Range("B5").Value = 1534
$ws.Range("D5").Value = 165
$ws.Range("E5").Value = 168
$ws.Range("H5").Value = 337
$ws.Range("I5").Value = 21607

$ws.Range("B6").Value = 1715
$ws.Range("D6").Value = 290
$ws.Range("E6").Value = 296
$ws.Range("H6").Value = 3138
$ws.Range("I6").Value = 56945
$ws.Range("J6").Value = 9

$ws.Range("B7").Value = 123
$ws.Range("D7").Value = 51
$ws.Range("E7").Value = 51
$ws.Range("I7").Value = 8824

$ws.Range("B8").Value = 1026
$ws.Range("D8").Value = 904
$ws.Range("E8").Value = 1026
$ws.Range("H8").Value = 623
$ws.Range("I8").Value = 13279
$ws.Range("J8").Value = 18
$ws.Range("K8").Value = 74
$ws.Range("L8").Value = 1034

$ws.Range("B9").Value = 262

$ws.Range("B10").Value = 1337
$ws.Range("D10").Value = 510
$ws.Range("E10").Value = 554
$ws.Range("F10").Value = 44
$ws.Range("H10").Value = 5448
$ws.Range("I10").Value = 15945

$ws.Range("B11").Value = 1026
$ws.Range("D11").Value = 256
$ws.Range("E11").Value = 286
$ws.Range("H11").Value = 2265
$ws.Range("I11").Value = 14148

$ws.Range("B12").Value = 1394

$ws.Range("B13").Value = 855
$ws.Range("D13").Value = 437
$ws.Range("E13").Value = 441
$ws.Range("H13").Value = 73
$ws.Range("I13").Value = 12739

$ws.Range("B14").Value = 264

$ws.Range("B15").Value = 437
$ws.Range("D15").Value = 148
$ws.Range("E15").Value = 156
$ws.Range("H15").Value = 1584
$ws.Range("I15").Value = 21196
$ws.Range("J15").Value = 5

$ws.Range("B16").Value = 183
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = 8
$ws.Range("H16").Value = 3750
$ws.Range("I16").Value = 2000

$ws.Range("B17").Value = 1771
$ws.Range("D17").Value = 35
$ws.Range("E17").Value = 36
$ws.Range("H17").Value = 1250
$ws.Range("I17").Value = 15238
$ws.Range("J17").Value = 1

$ws.Range("B18").Value = 4405
$ws.Range("D18").Value = 489
$ws.Range("E18").Value = 1272
$ws.Range("H18").Value = 2137
$ws.Range("I18").Value = 16510
$ws.Range("J18").Value = 7
$ws.Range("K18").Value = 693
$ws.Range("L18").Value = 20014
